# Trading update: 2026-02-17 04:08:48
# Appends the newly-opened MarketMaking trade (#25) as row 26 to both the
# "All Trades" and "MarketMaking" worksheets.

$wb = $excel.ActiveWorkbook

$newRowNumber = 26

# Values for the new trade row, keyed by column letter.
$tradeNo      = 25
$tradeDate    = "2026-02-17"
$tradeTime    = "04:08:43"
$strategy     = "MarketMaking"
$side         = "UP"
$entryPrice   = 0.57
$exitPrice    = ""                                   # still open -> blank
$status       = "OPEN"
$pnlPct       = 0
$pnlUsd       = 0
$capitalAfter = 100.2954564381429
$entrySlip    = 0
$exitSlip     = 0
$confidence   = 0.6
$entryReason  = "Normal spread capture: 19600 bps"
$exitReason   = ""                                   # still open -> blank
$durationMin  = 0

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item($newRowNumber, 1).Value = $tradeNo

    # Column B holds a plain "yyyy-mm-dd" looking string. Force text
    # formatting while assigning it so Excel doesn't auto-convert it into a
    # date serial number, then restore the default "Normal" style so the
    # cell keeps using the workbook's default (unstyled) formatting.
    $ws.Cells.Item($newRowNumber, 2).NumberFormat = "@"
    $ws.Cells.Item($newRowNumber, 2).Value = $tradeDate
    $ws.Cells.Item($newRowNumber, 2).Style = "Normal"

    $ws.Cells.Item($newRowNumber, 3).Value = $tradeTime
    $ws.Cells.Item($newRowNumber, 4).Value = $strategy
    $ws.Cells.Item($newRowNumber, 5).Value = $side
    $ws.Cells.Item($newRowNumber, 6).Value = $entryPrice

    # Exit Price: the trade is still OPEN, so this cell is blank but present.
    $ws.Cells.Item($newRowNumber, 7).NumberFormat = "General"
    $ws.Cells.Item($newRowNumber, 7).Style = "Normal"

    $ws.Cells.Item($newRowNumber, 8).Value = $status
    $ws.Cells.Item($newRowNumber, 9).Value = $pnlPct
    $ws.Cells.Item($newRowNumber, 10).Value = $pnlUsd
    $ws.Cells.Item($newRowNumber, 11).Value = $capitalAfter
    $ws.Cells.Item($newRowNumber, 12).Value = $entrySlip
    $ws.Cells.Item($newRowNumber, 13).Value = $exitSlip
    $ws.Cells.Item($newRowNumber, 14).Value = $confidence
    $ws.Cells.Item($newRowNumber, 15).Value = $entryReason

    # Exit Reason: blank but present, same treatment as Exit Price.
    $ws.Cells.Item($newRowNumber, 16).NumberFormat = "General"
    $ws.Cells.Item($newRowNumber, 16).Style = "Normal"

    $ws.Cells.Item($newRowNumber, 17).Value = $durationMin
}

Write-Output "Appended trade #$tradeNo to All Trades and MarketMaking sheets"
